# Fractal.Calc.xlsx - "refactor: code cleanup" / "refactor rename index"
#
# Renames the "index" column/header on the Fractal(2) sheet to "i" and
# re-bases the index values so the series starts at 0 instead of 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fractal(2)")

# Rename the header cell (also syncs the ListObject/table column name
# from "index" to "i").
$ws.Range("A1").Value = "i"

# Re-base the index column: subtract 1 from every data row (A2:A503),
# so row 2 becomes 0, row 3 becomes 1, ... row 503 becomes 501.
for ($r = 2; $r -le 503; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value() - 1
}

# Narrow column A now that it only holds 1-3 digit numbers.
$ws.Range("A:A").ColumnWidth = 3.14

# Move the active selection from I1 to M1.
$null = $ws.Range("M1").Select()
